# UI Run test based of Tags - Krishnaveni Vivekanandan
#
# Renames the auto-generated respondent first name from "Krish - Auto1" to
# "Krishnaveni - Auto1", tidies up the alt-contact first/last name test data
# ("John - Auto1" -> "John1", keeps "Updik1"), fixes the active-cell
# selection, widens the (now longer) respFirstName column, and removes a
# duplicate hyperlink that had been attached twice to M2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits -------------------------------------------------
# Order matters for how new entries land in the shared-string table:
# touch V2 before I2 so the newly introduced strings are appended in the
# same order as the reference workbook (John1 before Krishnaveni - Auto1).
$ws.Range("V2").Value = "John1"
$ws.Range("I2").Value = "Krishnaveni - Auto1"

# --- Selection ----------------------------------------------------------
$ws.Range("F9").Select()

# --- Column width (respFirstName got noticeably longer) -----------------
$ws.Columns.Item(9).ColumnWidth = 16.15

# --- Hyperlinks: drop the duplicate mailto link on M2 --------------------
# The engine's Hyperlink.Delete() on an individual item is a no-op, so
# rebuild the collection from scratch, keeping exactly one link per cell.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("Y2"), "mailto:John@gmail.com")
$ws.Hyperlinks.Add($ws.Range("AK2"), "mailto:JohnMS@gmail.com")
$ws.Hyperlinks.Add($ws.Range("M2"), "mailto:Akil1@gmail.com")
